$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "1) If you were to review the following code, ... Please"
$p1xml = '<w:p ' + $wNs + ' w:rsidR="00307527" w:rsidRPr="00307527" w:rsidRDefault="00307527" w:rsidP="00307527">' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="00307527">' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>1) If you were to review the following code, what feedback would you give? Please</w:t>' +
    '</w:r>' +
'</w:p>'
$d.Paragraphs.Item(1).Range.InsertXML($p1xml)

# --- Paragraph 2: "be specific and indicate any errors that would occur as well as other best practices"
$p2xml = '<w:p ' + $wNs + ' w:rsidR="00307527" w:rsidRPr="00307527" w:rsidRDefault="00307527" w:rsidP="00307527">' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="00307527">' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>be</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidRPr="00307527">' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t xml:space="preserve"> specific and indicate any errors that would occur as well as other best practices</w:t>' +
    '</w:r>' +
'</w:p>'
$d.Paragraphs.Item(2).Range.InsertXML($p2xml)

# --- Paragraph 3: "and code refactoring that should be done."
$p3xml = '<w:p ' + $wNs + ' w:rsidR="005342BB" w:rsidRPr="00307527" w:rsidRDefault="00307527" w:rsidP="00307527">' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="00307527">' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>and</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidRPr="00307527">' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t xml:space="preserve"> code refactoring that should be done.</w:t>' +
    '</w:r>' +
'</w:p>'
$d.Paragraphs.Item(3).Range.InsertXML($p3xml)

# --- Paragraph 4: empty paragraph, just pPr additions
$p4xml = '<w:p ' + $wNs + ' w:rsidR="00307527" w:rsidRPr="00307527" w:rsidRDefault="00307527" w:rsidP="00307527">' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
'</w:p>'
$d.Paragraphs.Item(4).Range.InsertXML($p4xml)

# --- Insert two brand-new empty paragraphs after paragraph 4 (pushes the
#     original trailing empty paragraph further down); then fill them in with
#     InsertXML. The original trailing paragraph (now at index 7) is left
#     completely untouched, per the diff.
$d.Paragraphs.Item(4).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(5).Range.InsertParagraphAfter() | Out-Null

# --- New Paragraph 5: "2) Rewrite the code based on the feedback you provided in question 1. Please"
$p5xml_new = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t xml:space="preserve">2) </w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>Rewrite</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t xml:space="preserve"> the code based on the feedback you provided in question 1. </w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>Please</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
'</w:p>'
$d.Paragraphs.Item(5).Range.InsertXML($p5xml_new)

# --- New Paragraph 6: "include unit tests on your code."
$p6xml_new = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>include</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t xml:space="preserve"> unit tests on your code.</w:t>' +
    '</w:r>' +
'</w:p>'
$d.Paragraphs.Item(6).Range.InsertXML($p6xml_new)

# --- Paragraph 7 (original, final trailing empty paragraph) is left
#     completely untouched — not present in the diff's changed lines.
